$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "Jan-15"
$ws.Cells.Item(2,2).Value = 1
$ws.Cells.Item(2,3).Value = 256252
$ws.Cells.Item(2,4).Value = 59478.78
$ws.Cells.Item(2,5).Value = 315730.78
$ws.Cells.Item(3,1).Value = "Feb-15"
$ws.Cells.Item(3,2).Value = 1
$ws.Cells.Item(3,3).Value = 258646.23
$ws.Cells.Item(3,4).Value = 64905.18
$ws.Cells.Item(3,5).Value = 323551.41
$ws.Cells.Item(4,1).Value = "Mar-15"
$ws.Cells.Item(4,2).Value = 1
$ws.Cells.Item(4,3).Value = 256908.04
$ws.Cells.Item(4,4).Value = 62641.01
$ws.Cells.Item(4,5).Value = 319549.05
$ws.Cells.Item(5,1).Value = "Mar-15"
$ws.Cells.Item(5,2).Value = 2
$ws.Cells.Item(5,3).Value = 160742.94
$ws.Cells.Item(5,4).Value = 38995.74
$ws.Cells.Item(5,5).Value = 199738.68
$ws.Cells.Item(6,1).Value = "Apr-15"
$ws.Cells.Item(6,2).Value = 1
$ws.Cells.Item(6,3).Value = 256056.24
$ws.Cells.Item(6,4).Value = 64213.55
$ws.Cells.Item(6,5).Value = 320269.79
$ws.Cells.Item(7,1).Value = "Apr-15"
$ws.Cells.Item(7,2).Value = 2
$ws.Cells.Item(7,3).Value = 168995.27
$ws.Cells.Item(7,4).Value = 39896.76
$ws.Cells.Item(7,5).Value = 208892.03
$ws.Cells.Item(8,1).Value = "May-15"
$ws.Cells.Item(8,2).Value = 1
$ws.Cells.Item(8,3).Value = 260284.16
$ws.Cells.Item(8,4).Value = 57556.92
$ws.Cells.Item(8,5).Value = 317841.08
$ws.Cells.Item(9,1).Value = "May-15"
$ws.Cells.Item(9,2).Value = 2
$ws.Cells.Item(9,3).Value = 167195.28
$ws.Cells.Item(9,4).Value = 40195.44
$ws.Cells.Item(9,5).Value = 207390.72
$ws.Cells.Item(10,1).Value = "Jun-15"
$ws.Cells.Item(10,2).Value = 1
$ws.Cells.Item(10,3).Value = 256414.89
$ws.Cells.Item(10,4).Value = 61346.3
$ws.Cells.Item(10,5).Value = 317761.19
$ws.Cells.Item(11,1).Value = "Jun-15"
$ws.Cells.Item(11,2).Value = 2
$ws.Cells.Item(11,3).Value = 170945.33
$ws.Cells.Item(11,4).Value = 37761.88
$ws.Cells.Item(11,5).Value = 208707.21
$ws.Cells.Item(12,1).Value = "Jul-15"
$ws.Cells.Item(12,2).Value = 1
$ws.Cells.Item(12,3).Value = 267660.8
$ws.Cells.Item(12,4).Value = 61584.04
$ws.Cells.Item(12,5).Value = 329244.84
$ws.Cells.Item(13,1).Value = "Jul-15"
$ws.Cells.Item(13,2).Value = 2
$ws.Cells.Item(13,3).Value = 178058.27
$ws.Cells.Item(13,4).Value = 36528.44
$ws.Cells.Item(13,5).Value = 214586.71
$ws.Cells.Item(14,1).Value = "Aug-15"
$ws.Cells.Item(14,2).Value = 1
$ws.Cells.Item(14,3).Value = 262457.27
$ws.Cells.Item(14,4).Value = 62696.12
$ws.Cells.Item(14,5).Value = 325153.39
$ws.Cells.Item(15,1).Value = "Aug-15"
$ws.Cells.Item(15,2).Value = 2
$ws.Cells.Item(15,3).Value = 175032.39
$ws.Cells.Item(15,4).Value = 41534.43
$ws.Cells.Item(15,5).Value = 216566.82
$ws.Cells.Item(16,1).Value = "Sep-15"
$ws.Cells.Item(16,2).Value = 1
$ws.Cells.Item(16,3).Value = 258928.32
$ws.Cells.Item(16,4).Value = 61834.99
$ws.Cells.Item(16,5).Value = 320763.31
$ws.Cells.Item(17,1).Value = "Sep-15"
$ws.Cells.Item(17,2).Value = 2
$ws.Cells.Item(17,3).Value = 175960.43
$ws.Cells.Item(17,4).Value = 48544.6
$ws.Cells.Item(17,5).Value = 224505.03
$ws.Cells.Item(18,1).Value = "Oct-15"
$ws.Cells.Item(18,2).Value = 1
$ws.Cells.Item(18,3).Value = 259209.61
$ws.Cells.Item(18,4).Value = 58787.85
$ws.Cells.Item(18,5).Value = 317997.46
$ws.Cells.Item(19,1).Value = "Oct-15"
$ws.Cells.Item(19,2).Value = 2
$ws.Cells.Item(19,3).Value = 183109.08
$ws.Cells.Item(19,4).Value = 42289.25
$ws.Cells.Item(19,5).Value = 225398.33
$ws.Cells.Item(20,1).Value = "Nov-15"
$ws.Cells.Item(20,2).Value = 1
$ws.Cells.Item(20,3).Value = 259236.71
$ws.Cells.Item(20,4).Value = 54222.69
$ws.Cells.Item(20,5).Value = 313459.4
$ws.Cells.Item(21,1).Value = "Nov-15"
$ws.Cells.Item(21,2).Value = 2
$ws.Cells.Item(21,3).Value = 174740.62
$ws.Cells.Item(21,4).Value = 43348.3
$ws.Cells.Item(21,5).Value = 218088.92
$ws.Cells.Item(22,1).Value = "Dec-15"
$ws.Cells.Item(22,2).Value = 1
$ws.Cells.Item(22,3).Value = 260436.51
$ws.Cells.Item(22,4).Value = 58130.84
$ws.Cells.Item(22,5).Value = 318567.35
$ws.Cells.Item(23,1).Value = "Dec-15"
$ws.Cells.Item(23,2).Value = 2
$ws.Cells.Item(23,3).Value = 175739.22
$ws.Cells.Item(23,4).Value = 42108.26
$ws.Cells.Item(23,5).Value = 217847.48
$ws.Cells.Item(24,1).Value = "Jan-16"
$ws.Cells.Item(24,2).Value = 1
$ws.Cells.Item(24,3).Value = 248265.01
$ws.Cells.Item(24,4).Value = 57011.33
$ws.Cells.Item(24,5).Value = 305276.34
$ws.Cells.Item(25,1).Value = "Jan-16"
$ws.Cells.Item(25,2).Value = 2
$ws.Cells.Item(25,3).Value = 173373.97
$ws.Cells.Item(25,4).Value = 45299.31
$ws.Cells.Item(25,5).Value = 218673.28
$ws.Cells.Item(26,1).Value = "Feb-16"
$ws.Cells.Item(26,2).Value = 1
$ws.Cells.Item(26,3).Value = 249659.04
$ws.Cells.Item(26,4).Value = 53004.35
$ws.Cells.Item(26,5).Value = 302663.39
$ws.Cells.Item(27,1).Value = "Feb-16"
$ws.Cells.Item(27,2).Value = 2
$ws.Cells.Item(27,3).Value = 174967.91
$ws.Cells.Item(27,4).Value = 42312.91
$ws.Cells.Item(27,5).Value = 217280.82
$ws.Cells.Item(28,1).Value = "Mar-16"
$ws.Cells.Item(28,2).Value = 1
$ws.Cells.Item(28,3).Value = 243872.28
$ws.Cells.Item(28,4).Value = 57653.59
$ws.Cells.Item(28,5).Value = 301525.87
$ws.Cells.Item(29,1).Value = "Mar-16"
$ws.Cells.Item(29,2).Value = 2
$ws.Cells.Item(29,3).Value = 178150
$ws.Cells.Item(29,4).Value = 44434.03
$ws.Cells.Item(29,5).Value = 222584.03
$ws.Cells.Item(30,1).Value = "Apr-16"
$ws.Cells.Item(30,2).Value = 1
$ws.Cells.Item(30,3).Value = 238906.38
$ws.Cells.Item(30,4).Value = 56976.27
$ws.Cells.Item(30,5).Value = 295882.65
$ws.Cells.Item(31,1).Value = "Apr-16"
$ws.Cells.Item(31,2).Value = 2
$ws.Cells.Item(31,3).Value = 175385.34
$ws.Cells.Item(31,4).Value = 45371.2
$ws.Cells.Item(31,5).Value = 220756.54
$ws.Cells.Item(32,1).Value = "May-16"
$ws.Cells.Item(32,2).Value = 1
$ws.Cells.Item(32,3).Value = 237134.48
$ws.Cells.Item(32,4).Value = 56748.89
$ws.Cells.Item(32,5).Value = 293883.37
$ws.Cells.Item(33,1).Value = "May-16"
$ws.Cells.Item(33,2).Value = 2
$ws.Cells.Item(33,3).Value = 166381.91
$ws.Cells.Item(33,4).Value = 39407.9
$ws.Cells.Item(33,5).Value = 205789.81
$ws.Cells.Item(34,1).Value = "Jun-16"
$ws.Cells.Item(34,2).Value = 1
$ws.Cells.Item(34,3).Value = 241764.6
$ws.Cells.Item(34,4).Value = 60886.09
$ws.Cells.Item(34,5).Value = 302650.69
$ws.Cells.Item(35,1).Value = "Jun-16"
$ws.Cells.Item(35,2).Value = 2
$ws.Cells.Item(35,3).Value = 166379.95
$ws.Cells.Item(35,4).Value = 42616.75
$ws.Cells.Item(35,5).Value = 208996.7
$ws.Cells.Item(36,1).Value = "Jul-16"
$ws.Cells.Item(36,2).Value = 1
$ws.Cells.Item(36,3).Value = 241689.72
$ws.Cells.Item(36,4).Value = 57185.97
$ws.Cells.Item(36,5).Value = 298875.69
$ws.Cells.Item(37,1).Value = "Jul-16"
$ws.Cells.Item(37,2).Value = 2
$ws.Cells.Item(37,3).Value = 168269.36
$ws.Cells.Item(37,4).Value = 45698.84
$ws.Cells.Item(37,5).Value = 213968.2
$ws.Cells.Item(38,1).Value = "Aug-16"
$ws.Cells.Item(38,2).Value = 1
$ws.Cells.Item(38,3).Value = 239990.58
$ws.Cells.Item(38,4).Value = 61696.1
$ws.Cells.Item(38,5).Value = 301686.68
$ws.Cells.Item(39,1).Value = "Aug-16"
$ws.Cells.Item(39,2).Value = 2
$ws.Cells.Item(39,3).Value = 155189.91
$ws.Cells.Item(39,4).Value = 45241.16
$ws.Cells.Item(39,5).Value = 200431.07
$ws.Cells.Item(40,1).Value = "Sep-16"
$ws.Cells.Item(40,2).Value = 1
$ws.Cells.Item(40,3).Value = 240313.86
$ws.Cells.Item(40,4).Value = 58766.4
$ws.Cells.Item(40,5).Value = 299080.26
$ws.Cells.Item(41,1).Value = "Sep-16"
$ws.Cells.Item(41,2).Value = 2
$ws.Cells.Item(41,3).Value = 161965.42
$ws.Cells.Item(41,4).Value = 37059.63
$ws.Cells.Item(41,5).Value = 199025.05
$ws.Cells.Item(42,1).Value = "Oct-16"
$ws.Cells.Item(42,2).Value = 1
$ws.Cells.Item(42,3).Value = 240012.32
$ws.Cells.Item(42,4).Value = 59436.24
$ws.Cells.Item(42,5).Value = 299448.56
$ws.Cells.Item(43,1).Value = "Oct-16"
$ws.Cells.Item(43,2).Value = 2
$ws.Cells.Item(43,3).Value = 172751.37
$ws.Cells.Item(43,4).Value = 37343.68
$ws.Cells.Item(43,5).Value = 210095.05
$ws.Cells.Item(44,1).Value = "Nov-16"
$ws.Cells.Item(44,2).Value = 1
$ws.Cells.Item(44,3).Value = 238319.46
$ws.Cells.Item(44,4).Value = 57591.61
$ws.Cells.Item(44,5).Value = 295911.07
$ws.Cells.Item(45,1).Value = "Nov-16"
$ws.Cells.Item(45,2).Value = 2
$ws.Cells.Item(45,3).Value = 176416.56
$ws.Cells.Item(45,4).Value = 40172.23
$ws.Cells.Item(45,5).Value = 216588.79
$ws.Cells.Item(46,1).Value = "Dec-16"
$ws.Cells.Item(46,2).Value = 1
$ws.Cells.Item(46,3).Value = 240254.26
$ws.Cells.Item(46,4).Value = 60005.85
$ws.Cells.Item(46,5).Value = 300260.11
$ws.Cells.Item(47,1).Value = "Dec-16"
$ws.Cells.Item(47,2).Value = 2
$ws.Cells.Item(47,3).Value = 175283.73
$ws.Cells.Item(47,4).Value = 43737.57
$ws.Cells.Item(47,5).Value = 219021.3
$ws.Cells.Item(48,1).Value = "Jan-17"
$ws.Cells.Item(48,2).Value = 1
$ws.Cells.Item(48,3).Value = 241911.7
$ws.Cells.Item(48,4).Value = 59090.7
$ws.Cells.Item(48,5).Value = 301002.4
$ws.Cells.Item(49,1).Value = "Jan-17"
$ws.Cells.Item(49,2).Value = 2
$ws.Cells.Item(49,3).Value = 173387.7
$ws.Cells.Item(49,4).Value = 38527.77
$ws.Cells.Item(49,5).Value = 211915.47
$ws.Cells.Item(50,1).Value = "Feb-17"
$ws.Cells.Item(50,2).Value = 1
$ws.Cells.Item(50,3).Value = 249348.41
$ws.Cells.Item(50,4).Value = 60792.92
$ws.Cells.Item(50,5).Value = 310141.33
$ws.Cells.Item(51,1).Value = "Feb-17"
$ws.Cells.Item(51,2).Value = 2
$ws.Cells.Item(51,3).Value = 175544.71
$ws.Cells.Item(51,4).Value = 39869.01
$ws.Cells.Item(51,5).Value = 215413.72
$ws.Cells.Item(52,1).Value = "Mar-17"
$ws.Cells.Item(52,2).Value = 1
$ws.Cells.Item(52,3).Value = 245877
$ws.Cells.Item(52,4).Value = 59932.55
$ws.Cells.Item(52,5).Value = 305809.55
$ws.Cells.Item(53,1).Value = "Mar-17"
$ws.Cells.Item(53,2).Value = 2
$ws.Cells.Item(53,3).Value = 175935.54
$ws.Cells.Item(53,4).Value = 39048.2
$ws.Cells.Item(53,5).Value = 214983.74
$ws.Cells.Item(54,1).Value = "Apr-17"
$ws.Cells.Item(54,2).Value = 1
$ws.Cells.Item(54,3).Value = 251825.11
$ws.Cells.Item(54,4).Value = 62354.46
$ws.Cells.Item(54,5).Value = 314179.57
$ws.Cells.Item(55,1).Value = "Apr-17"
$ws.Cells.Item(55,2).Value = 2
$ws.Cells.Item(55,3).Value = 188537.22
$ws.Cells.Item(55,4).Value = 41859.57
$ws.Cells.Item(55,5).Value = 230396.79
$ws.Cells.Item(56,1).Value = "May-17"
$ws.Cells.Item(56,2).Value = 1
$ws.Cells.Item(56,3).Value = 257984.21
$ws.Cells.Item(56,4).Value = 62318.97
$ws.Cells.Item(56,5).Value = 320303.18
$ws.Cells.Item(57,1).Value = "May-17"
$ws.Cells.Item(57,2).Value = 2
$ws.Cells.Item(57,3).Value = 175769.02
$ws.Cells.Item(57,4).Value = 42809.13
$ws.Cells.Item(57,5).Value = 218578.15
$ws.Cells.Item(58,1).Value = "Jun-17"
$ws.Cells.Item(58,2).Value = 1
$ws.Cells.Item(58,3).Value = 260625.04
$ws.Cells.Item(58,4).Value = 70422.25
$ws.Cells.Item(58,5).Value = 331047.29
$ws.Cells.Item(59,1).Value = "Jun-17"
$ws.Cells.Item(59,2).Value = 2
$ws.Cells.Item(59,3).Value = 178798.33
$ws.Cells.Item(59,4).Value = 44414.19
$ws.Cells.Item(59,5).Value = 223212.52
$ws.Cells.Item(60,1).Value = "Jul-17"
$ws.Cells.Item(60,2).Value = 1
$ws.Cells.Item(60,3).Value = 261607.26
$ws.Cells.Item(60,4).Value = 63084.26
$ws.Cells.Item(60,5).Value = 324691.52
$ws.Cells.Item(61,1).Value = "Jul-17"
$ws.Cells.Item(61,2).Value = 2
$ws.Cells.Item(61,3).Value = 187432.86
$ws.Cells.Item(61,4).Value = 46663.82
$ws.Cells.Item(61,5).Value = 234096.68
$ws.Cells.Item(62,1).Value = "Aug-17"
$ws.Cells.Item(62,2).Value = 1
$ws.Cells.Item(62,3).Value = 266230.15
$ws.Cells.Item(62,4).Value = 60189.49
$ws.Cells.Item(62,5).Value = 326419.64
$ws.Cells.Item(63,1).Value = "Aug-17"
$ws.Cells.Item(63,2).Value = 2
$ws.Cells.Item(63,3).Value = 192965.39
$ws.Cells.Item(63,4).Value = 40978.34
$ws.Cells.Item(63,5).Value = 233943.73
$ws.Cells.Item(64,1).Value = "Sep-17"
$ws.Cells.Item(64,2).Value = 1
$ws.Cells.Item(64,3).Value = 266774.3
$ws.Cells.Item(64,4).Value = 68294.15
$ws.Cells.Item(64,5).Value = 335068.45
$ws.Cells.Item(65,1).Value = "Sep-17"
$ws.Cells.Item(65,2).Value = 2
$ws.Cells.Item(65,3).Value = 182734.07
$ws.Cells.Item(65,4).Value = 40657.16
$ws.Cells.Item(65,5).Value = 223391.23
$ws.Cells.Item(66,1).Value = "Oct-17"
$ws.Cells.Item(66,2).Value = 1
$ws.Cells.Item(66,3).Value = 263510.02
$ws.Cells.Item(66,4).Value = 63688.93
$ws.Cells.Item(66,5).Value = 327198.95
$ws.Cells.Item(67,1).Value = "Oct-17"
$ws.Cells.Item(67,2).Value = 2
$ws.Cells.Item(67,3).Value = 174064.5
$ws.Cells.Item(67,4).Value = 42315.52
$ws.Cells.Item(67,5).Value = 216380.02
$ws.Cells.Item(68,1).Value = "Nov-17"
$ws.Cells.Item(68,2).Value = 1
$ws.Cells.Item(68,3).Value = 254558.82
$ws.Cells.Item(68,4).Value = 64080.31
$ws.Cells.Item(68,5).Value = 318639.13
$ws.Cells.Item(69,1).Value = "Nov-17"
$ws.Cells.Item(69,2).Value = 2
$ws.Cells.Item(69,3).Value = 172093.88
$ws.Cells.Item(69,4).Value = 47469.49
$ws.Cells.Item(69,5).Value = 219563.37
$ws.Cells.Item(70,1).Value = "Dec-17"
$ws.Cells.Item(70,2).Value = 1
$ws.Cells.Item(70,3).Value = 250056.43
$ws.Cells.Item(70,4).Value = 61480.02
$ws.Cells.Item(70,5).Value = 311536.45
$ws.Cells.Item(71,1).Value = "Dec-17"
$ws.Cells.Item(71,2).Value = 2
$ws.Cells.Item(71,3).Value = 178887.94
$ws.Cells.Item(71,4).Value = 41287.35
$ws.Cells.Item(71,5).Value = 220175.29
